$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")

$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1   # top
$c1.Borders.Item(9).LineStyle = 1   # top+bottom -> FINAL border4 idx2

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item(8).LineStyle = 1   # top only -> idx3 (border2)
$d1.Borders.Item(9).LineStyle = 1   # top+bottom -> should now match idx2 (border4, C1's state) and DISCARD idx3? let's see
